$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 2).Value = 5737536
$ws.Cells.Item(8, 6).Value = 'Juventus'
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 11).Value = 'A'
$ws.Cells.Item(8, 12).Value = 4
$ws.Cells.Item(8, 13).Value = 3.6
$ws.Cells.Item(8, 14).Value = 1.85
$ws.Cells.Item(8, 15).Value = 5.25
$ws.Cells.Item(8, 16).Value = 4.2
$ws.Cells.Item(8, 17).Value = 1.615
$ws.Cells.Item(8, 18).Value = 0.75
$ws.Cells.Item(8, 19).Value = 2.08
$ws.Cells.Item(8, 20).Value = 1.82
$ws.Cells.Item(8, 21).Value = 2.75
$ws.Cells.Item(8, 22).Value = 1.91
$ws.Cells.Item(8, 23).Value = 1.99
$ws.Cells.Item(8, 24).Value = -1
$ws.Cells.Item(8, 26).Value = 0.615
$ws.Cells.Item(8, 27).Value = -0.5
$ws.Cells.Item(8, 28).Value = 0.41
$ws.Cells.Item(8, 29).Value = -1
$ws.Cells.Item(8, 30).Value = 0.99
$ws.Cells.Item(9, 2).Value = 5713502
$ws.Cells.Item(9, 6).Value = 'Verona'
$ws.Cells.Item(9, 7).Value = 3
$ws.Cells.Item(9, 8).Value = 1
$ws.Cells.Item(9, 11).Value = 'H'
$ws.Cells.Item(9, 12).Value = 1.727
$ws.Cells.Item(9, 13).Value = 3.8
$ws.Cells.Item(9, 14).Value = 4.333
$ws.Cells.Item(9, 15).Value = 1.666
$ws.Cells.Item(9, 16).Value = 4
$ws.Cells.Item(9, 17).Value = 5
$ws.Cells.Item(9, 18).Value = -0.75
$ws.Cells.Item(9, 19).Value = 1.88
$ws.Cells.Item(9, 20).Value = 2.02
$ws.Cells.Item(9, 22).Value = 1.95
$ws.Cells.Item(9, 23).Value = 1.95
$ws.Cells.Item(9, 24).Value = 0.6659999999999999
$ws.Cells.Item(9, 26).Value = -1
$ws.Cells.Item(9, 27).Value = 0.8799999999999999
$ws.Cells.Item(9, 28).Value = -1
$ws.Cells.Item(9, 29).Value = 0.95
$ws.Cells.Item(10, 2).Value = 5706213
$ws.Cells.Item(10, 6).Value = 'Bologna'
$ws.Cells.Item(10, 7).Value = 2
$ws.Cells.Item(10, 8).Value = 3
$ws.Cells.Item(10, 9).Value = 1
$ws.Cells.Item(10, 11).Value = 'A'
$ws.Cells.Item(10, 12).Value = 2.625
$ws.Cells.Item(10, 13).Value = 3.2
$ws.Cells.Item(10, 14).Value = 2.6
$ws.Cells.Item(10, 15).Value = 2.45
$ws.Cells.Item(10, 16).Value = 3.5
$ws.Cells.Item(10, 17).Value = 2.75
$ws.Cells.Item(10, 18).Value = 0
$ws.Cells.Item(10, 19).Value = 1.87
$ws.Cells.Item(10, 20).Value = 2.03
$ws.Cells.Item(10, 21).Value = 2.5
$ws.Cells.Item(10, 22).Value = 2.07
$ws.Cells.Item(10, 23).Value = 1.83
$ws.Cells.Item(10, 24).Value = -1
$ws.Cells.Item(10, 26).Value = 1.75
$ws.Cells.Item(10, 27).Value = -1
$ws.Cells.Item(10, 28).Value = 1.03
$ws.Cells.Item(10, 29).Value = 1.07
$ws.Cells.Item(11, 2).Value = 5706211
$ws.Cells.Item(11, 6).Value = 'Monza'
$ws.Cells.Item(11, 7).Value = 5
$ws.Cells.Item(11, 8).Value = 2
$ws.Cells.Item(11, 9).Value = 2
$ws.Cells.Item(11, 11).Value = 'H'
$ws.Cells.Item(11, 12).Value = 1.5
$ws.Cells.Item(11, 13).Value = 4.5
$ws.Cells.Item(11, 14).Value = 5.5
$ws.Cells.Item(11, 15).Value = 1.5
$ws.Cells.Item(11, 16).Value = 5
$ws.Cells.Item(11, 17).Value = 5.75
$ws.Cells.Item(11, 18).Value = -1.25
$ws.Cells.Item(11, 19).Value = 2.05
$ws.Cells.Item(11, 20).Value = 1.85
$ws.Cells.Item(11, 21).Value = 3
$ws.Cells.Item(11, 22).Value = 1.88
$ws.Cells.Item(11, 23).Value = 2.02
$ws.Cells.Item(11, 24).Value = 0.5
$ws.Cells.Item(11, 26).Value = -1
$ws.Cells.Item(11, 27).Value = 1.05
$ws.Cells.Item(11, 28).Value = -1
$ws.Cells.Item(11, 29).Value = 0.8799999999999999
$ws.Cells.Item(11, 30).Value = -1
$ws.Cells.Item(16, 2).Value = 6868443
$ws.Cells.Item(16, 5).Value = 'Roma'
$ws.Cells.Item(16, 6).Value = 'Salernitana'
$ws.Cells.Item(16, 7).Value = 2
$ws.Cells.Item(16, 9).Value = 1
$ws.Cells.Item(16, 10).Value = 1
$ws.Cells.Item(16, 11).Value = 'D'
$ws.Cells.Item(16, 12).Value = 1.5
$ws.Cells.Item(16, 13).Value = 4
$ws.Cells.Item(16, 14).Value = 5.5
$ws.Cells.Item(16, 15).Value = 1.5
$ws.Cells.Item(16, 16).Value = 4
$ws.Cells.Item(16, 17).Value = 7.5
$ws.Cells.Item(16, 18).Value = -1
$ws.Cells.Item(16, 19).Value = 1.93
$ws.Cells.Item(16, 20).Value = 1.97
$ws.Cells.Item(16, 21).Value = 2.25
$ws.Cells.Item(16, 22).Value = 1.97
$ws.Cells.Item(16, 23).Value = 1.93
$ws.Cells.Item(16, 25).Value = 3
$ws.Cells.Item(16, 26).Value = -1
$ws.Cells.Item(16, 28).Value = 0.97
$ws.Cells.Item(16, 29).Value = 0.97
$ws.Cells.Item(16, 30).Value = -1
$ws.Cells.Item(17, 2).Value = 6868444
$ws.Cells.Item(17, 5).Value = 'Sassuolo'
$ws.Cells.Item(17, 6).Value = 'Udinese'
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 'A'
$ws.Cells.Item(17, 12).Value = 3
$ws.Cells.Item(17, 13).Value = 3.3
$ws.Cells.Item(17, 14).Value = 2.2
$ws.Cells.Item(17, 15).Value = 3.25
$ws.Cells.Item(17, 16).Value = 3.5
$ws.Cells.Item(17, 17).Value = 2.2
$ws.Cells.Item(17, 18).Value = 0.25
$ws.Cells.Item(17, 19).Value = 1.98
$ws.Cells.Item(17, 20).Value = 1.92
$ws.Cells.Item(17, 21).Value = 2.75
$ws.Cells.Item(17, 22).Value = 1.93
$ws.Cells.Item(17, 23).Value = 1.97
$ws.Cells.Item(17, 25).Value = -1
$ws.Cells.Item(17, 26).Value = 1.2
$ws.Cells.Item(17, 28).Value = 0.9199999999999999
$ws.Cells.Item(17, 29).Value = -1
$ws.Cells.Item(17, 30).Value = 0.97
$ws.Cells.Item(18, 5).Value = 'Atalanta'
$ws.Cells.Item(19, 5).Value = 'AC Milan'
$ws.Cells.Item(21, 6).Value = 'Lecce'
$ws.Cells.Item(23, 6).Value = 'Udinese'
$ws.Cells.Item(25, 5).Value = 'Lecce'
$ws.Cells.Item(26, 6).Value = 'Atalanta'
$ws.Cells.Item(30, 6).Value = 'AC Milan'
$ws.Cells.Item(33, 6).Value = 'Lecce'
$ws.Cells.Item(35, 5).Value = 'AC Milan'
$ws.Cells.Item(36, 5).Value = 'Udinese'
$ws.Cells.Item(40, 5).Value = 'Atalanta'
$ws.Cells.Item(43, 6).Value = 'Lecce'
$ws.Cells.Item(45, 6).Value = 'AC Milan'
$ws.Cells.Item(46, 6).Value = 'Atalanta'
$ws.Cells.Item(48, 6).Value = 'Udinese'
$ws.Cells.Item(53, 5).Value = 'Atalanta'
$ws.Cells.Item(54, 5).Value = 'Lecce'
$ws.Cells.Item(58, 5).Value = 'Udinese'
$ws.Cells.Item(59, 5).Value = 'AC Milan'
$ws.Cells.Item(62, 6).Value = 'Atalanta'
$ws.Cells.Item(63, 2).Value = 6868622
$ws.Cells.Item(63, 5).Value = 'Cagliari'
$ws.Cells.Item(63, 6).Value = 'Lecce'
$ws.Cells.Item(63, 8).Value = 3
$ws.Cells.Item(63, 10).Value = 2
$ws.Cells.Item(63, 11).Value = 'A'
$ws.Cells.Item(63, 12).Value = 6
$ws.Cells.Item(63, 13).Value = 4
$ws.Cells.Item(63, 14).Value = 1.533
$ws.Cells.Item(63, 15).Value = 4.333
$ws.Cells.Item(63, 16).Value = 3.5
$ws.Cells.Item(63, 17).Value = 1.85
$ws.Cells.Item(63, 18).Value = 0.5
$ws.Cells.Item(63, 19).Value = 1.98
$ws.Cells.Item(63, 20).Value = 1.92
$ws.Cells.Item(63, 22).Value = 2.03
$ws.Cells.Item(63, 23).Value = 1.87
$ws.Cells.Item(63, 24).Value = -1
$ws.Cells.Item(63, 26).Value = 0.8500000000000001
$ws.Cells.Item(63, 27).Value = -1
$ws.Cells.Item(63, 28).Value = 0.9199999999999999
$ws.Cells.Item(63, 29).Value = 1.03
$ws.Cells.Item(63, 30).Value = -1
$ws.Cells.Item(64, 2).Value = 6868623
$ws.Cells.Item(64, 5).Value = 'Empoli'
$ws.Cells.Item(64, 6).Value = 'Salernitana'
$ws.Cells.Item(64, 7).Value = 1
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 1
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 'H'
$ws.Cells.Item(64, 12).Value = 2.25
$ws.Cells.Item(64, 13).Value = 3.1
$ws.Cells.Item(64, 14).Value = 3.4
$ws.Cells.Item(64, 15).Value = 2.7
$ws.Cells.Item(64, 16).Value = 3.2
$ws.Cells.Item(64, 17).Value = 2.75
$ws.Cells.Item(64, 18).Value = 0
$ws.Cells.Item(64, 19).Value = 1.91
$ws.Cells.Item(64, 20).Value = 1.99
$ws.Cells.Item(64, 21).Value = 2.25
$ws.Cells.Item(64, 22).Value = 1.88
$ws.Cells.Item(64, 23).Value = 2.02
$ws.Cells.Item(64, 24).Value = 1.7
$ws.Cells.Item(64, 26).Value = -1
$ws.Cells.Item(64, 27).Value = 0.9099999999999999
$ws.Cells.Item(64, 28).Value = -1
$ws.Cells.Item(64, 30).Value = 1.02
$ws.Cells.Item(65, 2).Value = 6868626
$ws.Cells.Item(65, 5).Value = 'Verona'
$ws.Cells.Item(65, 6).Value = 'Udinese'
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 1
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 1
$ws.Cells.Item(65, 12).Value = 4.7
$ws.Cells.Item(65, 14).Value = 1.666
$ws.Cells.Item(65, 15).Value = 4.2
$ws.Cells.Item(65, 21).Value = 2.5
$ws.Cells.Item(65, 22).Value = 2.06
$ws.Cells.Item(65, 23).Value = 1.84
$ws.Cells.Item(65, 29).Value = -1
$ws.Cells.Item(65, 30).Value = 0.8400000000000001
$ws.Cells.Item(67, 2).Value = 6868631
$ws.Cells.Item(67, 5).Value = 'Napoli'
$ws.Cells.Item(67, 6).Value = 'AC Milan'
$ws.Cells.Item(67, 7).Value = 4
$ws.Cells.Item(67, 8).Value = 1
$ws.Cells.Item(67, 9).Value = 2
$ws.Cells.Item(67, 12).Value = 1.333
$ws.Cells.Item(67, 13).Value = 5
$ws.Cells.Item(67, 14).Value = 8.5
$ws.Cells.Item(67, 15).Value = 1.4
$ws.Cells.Item(67, 16).Value = 5
$ws.Cells.Item(67, 17).Value = 7
$ws.Cells.Item(67, 18).Value = -1.25
$ws.Cells.Item(67, 19).Value = 1.93
$ws.Cells.Item(67, 20).Value = 1.97
$ws.Cells.Item(67, 21).Value = 3
$ws.Cells.Item(67, 22).Value = 1.99
$ws.Cells.Item(67, 23).Value = 1.91
$ws.Cells.Item(67, 24).Value = 0.3999999999999999
$ws.Cells.Item(67, 27).Value = 0.9299999999999999
$ws.Cells.Item(67, 29).Value = 0.99
$ws.Cells.Item(67, 30).Value = -1
$ws.Cells.Item(68, 2).Value = 6868629
$ws.Cells.Item(68, 5).Value = 'Lazio'
$ws.Cells.Item(68, 6).Value = 'Torino'
$ws.Cells.Item(68, 7).Value = 2
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 12).Value = 1.833
$ws.Cells.Item(68, 13).Value = 3.5
$ws.Cells.Item(68, 14).Value = 4.333
$ws.Cells.Item(68, 15).Value = 2.25
$ws.Cells.Item(68, 16).Value = 3.2
$ws.Cells.Item(68, 17).Value = 3.4
$ws.Cells.Item(68, 18).Value = -0.25
$ws.Cells.Item(68, 19).Value = 1.97
$ws.Cells.Item(68, 20).Value = 1.93
$ws.Cells.Item(68, 21).Value = 2.25
$ws.Cells.Item(68, 22).Value = 2.05
$ws.Cells.Item(68, 23).Value = 1.85
$ws.Cells.Item(68, 24).Value = 1.25
$ws.Cells.Item(68, 27).Value = 0.97
$ws.Cells.Item(68, 29).Value = -0.5
$ws.Cells.Item(68, 30).Value = 0.425
$ws.Cells.Item(71, 5).Value = 'Atalanta'
$ws.Cells.Item(72, 5).Value = 'Lecce'
$ws.Cells.Item(75, 5).Value = 'AC Milan'
$ws.Cells.Item(76, 5).Value = 'Udinese'
$ws.Cells.Item(81, 6).Value = 'AC Milan'
$ws.Cells.Item(82, 5).Value = 'Atalanta'
$ws.Cells.Item(85, 6).Value = 'Lecce'
$ws.Cells.Item(87, 6).Value = 'Udinese'
$ws.Cells.Item(97, 5).Value = 'Udinese'
$ws.Cells.Item(98, 5).Value = 'Lecce'
$ws.Cells.Item(99, 5).Value = 'AC Milan'
$ws.Cells.Item(99, 6).Value = 'Atalanta'
$ws.Cells.Item(103, 5).Value = 'Atalanta'
$ws.Cells.Item(106, 6).Value = 'AC Milan'
$ws.Cells.Item(108, 6).Value = 'Lecce'
$ws.Cells.Item(109, 6).Value = 'Udinese'
$ws.Cells.Item(113, 5).Value = 'Udinese'
$ws.Cells.Item(114, 5).Value = 'Lecce'
$ws.Cells.Item(114, 6).Value = 'AC Milan'
$ws.Cells.Item(117, 6).Value = 'Atalanta'
$ws.Cells.Item(123, 5).Value = 'Atalanta'
$ws.Cells.Item(123, 6).Value = 'Lecce'
$ws.Cells.Item(127, 5).Value = 'AC Milan'
$ws.Cells.Item(127, 6).Value = 'Udinese'
$ws.Cells.Item(132, 5).Value = 'Udinese'
$ws.Cells.Item(133, 5).Value = 'Lecce'
$ws.Cells.Item(137, 6).Value = 'AC Milan'
$ws.Cells.Item(139, 6).Value = 'Atalanta'
$ws.Cells.Item(144, 5).Value = 'Lecce'
$ws.Cells.Item(145, 5).Value = 'Atalanta'
$ws.Cells.Item(146, 2).Value = 6868716
$ws.Cells.Item(146, 5).Value = 'AC Milan'
$ws.Cells.Item(146, 6).Value = 'Verona'
$ws.Cells.Item(146, 8).Value = 3
$ws.Cells.Item(146, 10).Value = 1
$ws.Cells.Item(146, 11).Value = 'D'
$ws.Cells.Item(146, 12).Value = 1.909
$ws.Cells.Item(146, 13).Value = 3.3
$ws.Cells.Item(146, 14).Value = 4.25
$ws.Cells.Item(146, 15).Value = 1.8
$ws.Cells.Item(146, 16).Value = 3.6
$ws.Cells.Item(146, 17).Value = 4.5
$ws.Cells.Item(146, 18).Value = -0.75
$ws.Cells.Item(146, 19).Value = 2.07
$ws.Cells.Item(146, 20).Value = 1.83
$ws.Cells.Item(146, 21).Value = 2.5
$ws.Cells.Item(146, 22).Value = 2.04
$ws.Cells.Item(146, 23).Value = 1.86
$ws.Cells.Item(146, 24).Value = -1
$ws.Cells.Item(146, 25).Value = 2.6
$ws.Cells.Item(146, 27).Value = -1
$ws.Cells.Item(146, 28).Value = 0.8300000000000001
$ws.Cells.Item(146, 29).Value = 1.04
$ws.Cells.Item(146, 30).Value = -1
$ws.Cells.Item(147, 2).Value = 6868707
$ws.Cells.Item(147, 5).Value = 'Fiorentina'
$ws.Cells.Item(147, 6).Value = 'Salernitana'
$ws.Cells.Item(147, 8).Value = 0
$ws.Cells.Item(147, 10).Value = 0
$ws.Cells.Item(147, 11).Value = 'H'
$ws.Cells.Item(147, 12).Value = 1.444
$ws.Cells.Item(147, 13).Value = 4.5
$ws.Cells.Item(147, 14).Value = 6.5
$ws.Cells.Item(147, 15).Value = 1.4
$ws.Cells.Item(147, 16).Value = 4.75
$ws.Cells.Item(147, 17).Value = 7.5
$ws.Cells.Item(147, 18).Value = -1.25
$ws.Cells.Item(147, 19).Value = 1.95
$ws.Cells.Item(147, 20).Value = 1.95
$ws.Cells.Item(147, 21).Value = 2.75
$ws.Cells.Item(147, 22).Value = 2.02
$ws.Cells.Item(147, 23).Value = 1.88
$ws.Cells.Item(147, 24).Value = 0.3999999999999999
$ws.Cells.Item(147, 25).Value = -1
$ws.Cells.Item(147, 27).Value = 0.95
$ws.Cells.Item(147, 28).Value = -1
$ws.Cells.Item(147, 29).Value = 0.51
$ws.Cells.Item(147, 30).Value = -0.5
$ws.Cells.Item(150, 6).Value = 'Udinese'
$ws.Cells.Item(153, 5).Value = 'Udinese'
$ws.Cells.Item(153, 6).Value = 'Lecce'
$ws.Cells.Item(154, 6).Value = 'AC Milan'
$ws.Cells.Item(159, 6).Value = 'Atalanta'
$ws.Cells.Item(162, 5).Value = 'Atalanta'
$ws.Cells.Item(165, 5).Value = 'Lecce'
$ws.Cells.Item(166, 2).Value = 6868729
$ws.Cells.Item(166, 5).Value = 'Fiorentina'
$ws.Cells.Item(166, 6).Value = 'Verona'
$ws.Cells.Item(166, 7).Value = 1
$ws.Cells.Item(166, 8).Value = 0
$ws.Cells.Item(166, 9).Value = 0
$ws.Cells.Item(166, 11).Value = 'H'
$ws.Cells.Item(166, 12).Value = 1.5
$ws.Cells.Item(166, 13).Value = 4.5
$ws.Cells.Item(166, 14).Value = 5.5
$ws.Cells.Item(166, 15).Value = 1.45
$ws.Cells.Item(166, 16).Value = 4.75
$ws.Cells.Item(166, 17).Value = 6
$ws.Cells.Item(166, 18).Value = -1
$ws.Cells.Item(166, 19).Value = 1.82
$ws.Cells.Item(166, 20).Value = 2.08
$ws.Cells.Item(166, 22).Value = 1.89
$ws.Cells.Item(166, 23).Value = 2.01
$ws.Cells.Item(166, 24).Value = 0.45
$ws.Cells.Item(166, 25).Value = -1
$ws.Cells.Item(166, 29).Value = -1
$ws.Cells.Item(166, 30).Value = 1.01
$ws.Cells.Item(167, 2).Value = 6868736
$ws.Cells.Item(167, 5).Value = 'AC Milan'
$ws.Cells.Item(167, 6).Value = 'Sassuolo'
$ws.Cells.Item(167, 7).Value = 2
$ws.Cells.Item(167, 8).Value = 2
$ws.Cells.Item(167, 9).Value = 1
$ws.Cells.Item(167, 11).Value = 'D'
$ws.Cells.Item(167, 12).Value = 2.15
$ws.Cells.Item(167, 13).Value = 3.2
$ws.Cells.Item(167, 14).Value = 3.3
$ws.Cells.Item(167, 15).Value = 2.5
$ws.Cells.Item(167, 16).Value = 3.3
$ws.Cells.Item(167, 17).Value = 2.8
$ws.Cells.Item(167, 18).Value = 0
$ws.Cells.Item(167, 19).Value = 1.83
$ws.Cells.Item(167, 20).Value = 2.07
$ws.Cells.Item(167, 22).Value = 1.92
$ws.Cells.Item(167, 23).Value = 1.98
$ws.Cells.Item(167, 24).Value = -1
$ws.Cells.Item(167, 25).Value = 2.3
$ws.Cells.Item(167, 29).Value = 0.9199999999999999
$ws.Cells.Item(167, 30).Value = -1
$ws.Cells.Item(170, 5).Value = 'Udinese'
$ws.Cells.Item(173, 2).Value = 6868742
$ws.Cells.Item(173, 5).Value = 'Monza'
$ws.Cells.Item(173, 6).Value = 'Fiorentina'
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 1
$ws.Cells.Item(173, 9).Value = 0
$ws.Cells.Item(173, 11).Value = 'A'
$ws.Cells.Item(173, 12).Value = 2.8
$ws.Cells.Item(173, 13).Value = 3.5
$ws.Cells.Item(173, 14).Value = 2.375
$ws.Cells.Item(173, 15).Value = 2.8
$ws.Cells.Item(173, 16).Value = 3.25
$ws.Cells.Item(173, 17).Value = 2.6
$ws.Cells.Item(173, 18).Value = 0
$ws.Cells.Item(173, 19).Value = 2.04
$ws.Cells.Item(173, 20).Value = 1.86
$ws.Cells.Item(173, 21).Value = 2.5
$ws.Cells.Item(173, 22).Value = 2.025
$ws.Cells.Item(173, 23).Value = 1.775
$ws.Cells.Item(173, 25).Value = -1
$ws.Cells.Item(173, 26).Value = 1.6
$ws.Cells.Item(173, 27).Value = -1
$ws.Cells.Item(173, 28).Value = 0.8600000000000001
$ws.Cells.Item(173, 29).Value = -1
$ws.Cells.Item(173, 30).Value = 0.7749999999999999
$ws.Cells.Item(174, 2).Value = 6868744
$ws.Cells.Item(174, 5).Value = 'Salernitana'
$ws.Cells.Item(174, 6).Value = 'Lecce'
$ws.Cells.Item(174, 7).Value = 2
$ws.Cells.Item(174, 8).Value = 2
$ws.Cells.Item(174, 9).Value = 1
$ws.Cells.Item(174, 11).Value = 'D'
$ws.Cells.Item(174, 12).Value = 5.75
$ws.Cells.Item(174, 13).Value = 4
$ws.Cells.Item(174, 14).Value = 1.533
$ws.Cells.Item(174, 15).Value = 6.5
$ws.Cells.Item(174, 16).Value = 4
$ws.Cells.Item(174, 17).Value = 1.5
$ws.Cells.Item(174, 18).Value = 1
$ws.Cells.Item(174, 19).Value = 2.05
$ws.Cells.Item(174, 20).Value = 1.85
$ws.Cells.Item(174, 21).Value = 2.75
$ws.Cells.Item(174, 22).Value = 1.89
$ws.Cells.Item(174, 23).Value = 2.01
$ws.Cells.Item(174, 25).Value = 3
$ws.Cells.Item(174, 26).Value = -1
$ws.Cells.Item(174, 27).Value = 1.05
$ws.Cells.Item(174, 28).Value = -1
$ws.Cells.Item(174, 29).Value = 0.8899999999999999
$ws.Cells.Item(174, 30).Value = -1
$ws.Cells.Item(176, 2).Value = 6868737
$ws.Cells.Item(176, 5).Value = 'Bologna'
$ws.Cells.Item(176, 8).Value = 0
$ws.Cells.Item(176, 11).Value = 'H'
$ws.Cells.Item(176, 12).Value = 2.7
$ws.Cells.Item(176, 13).Value = 3.25
$ws.Cells.Item(176, 14).Value = 2.45
$ws.Cells.Item(176, 15).Value = 2.8
$ws.Cells.Item(176, 16).Value = 3
$ws.Cells.Item(176, 17).Value = 2.75
$ws.Cells.Item(176, 18).Value = 0
$ws.Cells.Item(176, 19).Value = 1.95
$ws.Cells.Item(176, 20).Value = 1.95
$ws.Cells.Item(176, 21).Value = 2
$ws.Cells.Item(176, 22).Value = 1.85
$ws.Cells.Item(176, 23).Value = 2.05
$ws.Cells.Item(176, 24).Value = 1.8
$ws.Cells.Item(176, 25).Value = -1
$ws.Cells.Item(176, 27).Value = 0.95
$ws.Cells.Item(176, 28).Value = -1
$ws.Cells.Item(176, 29).Value = -1
$ws.Cells.Item(176, 30).Value = 1.05
$ws.Cells.Item(177, 2).Value = 6868746
$ws.Cells.Item(177, 5).Value = 'Torino'
$ws.Cells.Item(177, 6).Value = 'AC Milan'
$ws.Cells.Item(177, 8).Value = 1
$ws.Cells.Item(177, 11).Value = 'D'
$ws.Cells.Item(177, 12).Value = 1.85
$ws.Cells.Item(177, 13).Value = 3.3
$ws.Cells.Item(177, 14).Value = 4.2
$ws.Cells.Item(177, 15).Value = 1.65
$ws.Cells.Item(177, 16).Value = 3.75
$ws.Cells.Item(177, 17).Value = 5.5
$ws.Cells.Item(177, 18).Value = -0.75
$ws.Cells.Item(177, 19).Value = 1.85
$ws.Cells.Item(177, 20).Value = 2.05
$ws.Cells.Item(177, 21).Value = 2.25
$ws.Cells.Item(177, 22).Value = 1.99
$ws.Cells.Item(177, 23).Value = 1.91
$ws.Cells.Item(177, 24).Value = -1
$ws.Cells.Item(177, 25).Value = 2.75
$ws.Cells.Item(177, 27).Value = -1
$ws.Cells.Item(177, 28).Value = 1.05
$ws.Cells.Item(177, 29).Value = -0.5
$ws.Cells.Item(177, 30).Value = 0.455
$ws.Cells.Item(178, 6).Value = 'Atalanta'
$ws.Cells.Item(185, 5).Value = 'Udinese'
$ws.Cells.Item(185, 6).Value = 'Atalanta'
$ws.Cells.Item(186, 5).Value = 'AC Milan'
$ws.Cells.Item(188, 5).Value = 'Lecce'
$ws.Cells.Item(194, 5).Value = 'Atalanta'
$ws.Cells.Item(196, 6).Value = 'Lecce'
$ws.Cells.Item(198, 5).Value = 'AC Milan'
$ws.Cells.Item(200, 6).Value = 'Udinese'
$ws.Cells.Item(205, 6).Value = 'Atalanta'
$ws.Cells.Item(207, 6).Value = 'AC Milan'
$ws.Cells.Item(208, 5).Value = 'Lecce'
$ws.Cells.Item(209, 5).Value = 'Udinese'
$ws.Cells.Item(212, 5).Value = 'AC Milan'
$ws.Cells.Item(212, 6).Value = 'Lecce'
$ws.Cells.Item(216, 5).Value = 'Atalanta'
$ws.Cells.Item(218, 5).Value = 'Udinese'
$ws.Cells.Item(218, 6).Value = 'AC Milan'
$ws.Cells.Item(220, 5).Value = 'Lecce'
$ws.Cells.Item(221, 6).Value = 'Atalanta'
$ws.Cells.Item(227, 5).Value = 'Atalanta'
$ws.Cells.Item(229, 5).Value = 'AC Milan'
$ws.Cells.Item(230, 6).Value = 'Lecce'
$ws.Cells.Item(234, 5).Value = 'Udinese'
$ws.Cells.Item(243, 6).Value = 'Atalanta'
$ws.Cells.Item(244, 6).Value = 'Udinese'
$ws.Cells.Item(245, 5).Value = 'Lecce'
$ws.Cells.Item(246, 6).Value = 'AC Milan'
$ws.Cells.Item(248, 6).Value = 'Atalanta'
$ws.Cells.Item(252, 5).Value = 'Udinese'
$ws.Cells.Item(255, 5).Value = 'AC Milan'
$ws.Cells.Item(257, 6).Value = 'Lecce'
$ws.Cells.Item(262, 6).Value = 'AC Milan'
$ws.Cells.Item(265, 5).Value = 'Atalanta'
$ws.Cells.Item(266, 5).Value = 'Lecce'
$ws.Cells.Item(266, 6).Value = 'Udinese'
$ws.Cells.Item(270, 6).Value = 'Udinese'
$ws.Cells.Item(271, 6).Value = 'Lecce'
$ws.Cells.Item(272, 5).Value = 'AC Milan'
$ws.Cells.Item(277, 6).Value = 'Atalanta'
$ws.Cells.Item(278, 5).Value = 'Udinese'
$ws.Cells.Item(286, 5).Value = 'Atalanta'
$ws.Cells.Item(287, 5).Value = 'Lecce'
$ws.Cells.Item(288, 6).Value = 'Udinese'
$ws.Cells.Item(290, 6).Value = 'AC Milan'
$ws.Cells.Item(293, 5).Value = 'AC Milan'
$ws.Cells.Item(294, 6).Value = 'Atalanta'
$ws.Cells.Item(297, 6).Value = 'Lecce'
$ws.Cells.Item(300, 6).Value = 'Udinese'
$ws.Cells.Item(304, 6).Value = 'Lecce'
$ws.Cells.Item(307, 6).Value = 'AC Milan'
$ws.Cells.Item(308, 5).Value = 'Atalanta'
$ws.Cells.Item(311, 5).Value = 'Lecce'
$ws.Cells.Item(311, 6).Value = 'Atalanta'
$ws.Cells.Item(316, 6).Value = 'Udinese'
$ws.Cells.Item(319, 5).Value = 'AC Milan'
$ws.Cells.Item(321, 5).Value = 'Atalanta'
$ws.Cells.Item(325, 6).Value = 'Lecce'
$ws.Cells.Item(326, 5).Value = 'AC Milan'
$ws.Cells.Item(329, 5).Value = 'Udinese'
$ws.Cells.Item(333, 6).Value = 'AC Milan'
$ws.Cells.Item(334, 6).Value = 'Atalanta'
$ws.Cells.Item(337, 6).Value = 'Udinese'
$ws.Cells.Item(339, 5).Value = 'Lecce'
$ws.Cells.Item(341, 5).Value = 'Atalanta'
$ws.Cells.Item(342, 6).Value = 'Lecce'
$ws.Cells.Item(345, 6).Value = 'AC Milan'
$ws.Cells.Item(346, 5).Value = 'Udinese'
$ws.Cells.Item(353, 6).Value = 'Atalanta'
$ws.Cells.Item(356, 5).Value = 'Lecce'
$ws.Cells.Item(358, 6).Value = 'Udinese'
$ws.Cells.Item(359, 5).Value = 'AC Milan'
$ws.Cells.Item(362, 5).Value = 'Lecce'
$ws.Cells.Item(364, 2).Value = 7288907
$ws.Cells.Item(364, 5).Value = 'Verona'
$ws.Cells.Item(364, 6).Value = 'Torino'
$ws.Cells.Item(364, 7).Value = 1
$ws.Cells.Item(364, 8).Value = 2
$ws.Cells.Item(364, 10).Value = 0
$ws.Cells.Item(364, 11).Value = 'A'
$ws.Cells.Item(364, 12).Value = 3
$ws.Cells.Item(364, 13).Value = 3
$ws.Cells.Item(364, 14).Value = 2.5
$ws.Cells.Item(364, 15).Value = 3.1
$ws.Cells.Item(364, 16).Value = 2.9
$ws.Cells.Item(364, 17).Value = 2.55
$ws.Cells.Item(364, 18).Value = 0
$ws.Cells.Item(364, 19).Value = 2.15
$ws.Cells.Item(364, 20).Value = 1.79
$ws.Cells.Item(364, 21).Value = 2
$ws.Cells.Item(364, 22).Value = 2.02
$ws.Cells.Item(364, 23).Value = 1.88
$ws.Cells.Item(364, 24).Value = -1
$ws.Cells.Item(364, 26).Value = 1.55
$ws.Cells.Item(364, 27).Value = -1
$ws.Cells.Item(364, 28).Value = 0.79
$ws.Cells.Item(364, 29).Value = 1.02
$ws.Cells.Item(365, 2).Value = 7284809
$ws.Cells.Item(365, 5).Value = 'Genoa'
$ws.Cells.Item(365, 6).Value = 'Sassuolo'
$ws.Cells.Item(365, 7).Value = 2
$ws.Cells.Item(365, 8).Value = 1
$ws.Cells.Item(365, 10).Value = 1
$ws.Cells.Item(365, 11).Value = 'H'
$ws.Cells.Item(365, 12).Value = 2
$ws.Cells.Item(365, 13).Value = 3.5
$ws.Cells.Item(365, 14).Value = 3.6
$ws.Cells.Item(365, 15).Value = 2.25
$ws.Cells.Item(365, 16).Value = 3.4
$ws.Cells.Item(365, 17).Value = 3.2
$ws.Cells.Item(365, 18).Value = -0.25
$ws.Cells.Item(365, 19).Value = 2
$ws.Cells.Item(365, 20).Value = 1.9
$ws.Cells.Item(365, 21).Value = 2.5
$ws.Cells.Item(365, 22).Value = 2.025
$ws.Cells.Item(365, 23).Value = 1.775
$ws.Cells.Item(365, 24).Value = 1.25
$ws.Cells.Item(365, 26).Value = -1
$ws.Cells.Item(365, 27).Value = 1
$ws.Cells.Item(365, 28).Value = -1
$ws.Cells.Item(365, 29).Value = 1.025
$ws.Cells.Item(367, 5).Value = 'Udinese'
$ws.Cells.Item(368, 5).Value = 'Atalanta'
$ws.Cells.Item(368, 6).Value = 'AC Milan'
$ws.Cells.Item(371, 5).Value = 'Atalanta'
$ws.Cells.Item(371, 6).Value = 'Udinese'
$ws.Cells.Item(372, 6).Value = 'Lecce'
$ws.Cells.Item(374, 2).Value = 7284017
$ws.Cells.Item(374, 5).Value = 'AC Milan'
$ws.Cells.Item(374, 6).Value = 'Empoli'
$ws.Cells.Item(374, 7).Value = 1
$ws.Cells.Item(374, 10).Value = 0
$ws.Cells.Item(374, 11).Value = 'D'
$ws.Cells.Item(374, 12).Value = 2.05
$ws.Cells.Item(374, 13).Value = 3.3
$ws.Cells.Item(374, 14).Value = 3.6
$ws.Cells.Item(374, 15).Value = 1.909
$ws.Cells.Item(374, 16).Value = 3.3
$ws.Cells.Item(374, 17).Value = 4.333
$ws.Cells.Item(374, 18).Value = -0.5
$ws.Cells.Item(374, 19).Value = 1.98
$ws.Cells.Item(374, 20).Value = 1.92
$ws.Cells.Item(374, 21).Value = 2
$ws.Cells.Item(374, 22).Value = 1.85
$ws.Cells.Item(374, 23).Value = 2.05
$ws.Cells.Item(374, 25).Value = 2.3
$ws.Cells.Item(374, 26).Value = -1
$ws.Cells.Item(374, 28).Value = 0.9199999999999999
$ws.Cells.Item(374, 29).Value = 0
$ws.Cells.Item(374, 30).Value = 0
$ws.Cells.Item(375, 2).Value = 7288908
$ws.Cells.Item(375, 5).Value = 'Monza'
$ws.Cells.Item(375, 6).Value = 'Frosinone'
$ws.Cells.Item(375, 7).Value = 0
$ws.Cells.Item(375, 10).Value = 1
$ws.Cells.Item(375, 11).Value = 'A'
$ws.Cells.Item(375, 12).Value = 2.7
$ws.Cells.Item(375, 13).Value = 3.4
$ws.Cells.Item(375, 14).Value = 2.5
$ws.Cells.Item(375, 15).Value = 3
$ws.Cells.Item(375, 16).Value = 3.5
$ws.Cells.Item(375, 17).Value = 2.3
$ws.Cells.Item(375, 18).Value = 0.25
$ws.Cells.Item(375, 19).Value = 1.88
$ws.Cells.Item(375, 20).Value = 2.05
$ws.Cells.Item(375, 21).Value = 2.75
$ws.Cells.Item(375, 22).Value = 1.87
$ws.Cells.Item(375, 23).Value = 2.03
$ws.Cells.Item(375, 25).Value = -1
$ws.Cells.Item(375, 26).Value = 1.3
$ws.Cells.Item(375, 28).Value = 1.05
$ws.Cells.Item(375, 29).Value = -1
$ws.Cells.Item(375, 30).Value = 1.03
$ws.Cells.Item(383, 5).Value = 'Lecce'
$ws.Cells.Item(384, 5).Value = 'Udinese'
$ws.Cells.Item(385, 6).Value = 'Atalanta'
$ws.Cells.Item(386, 6).Value = 'AC Milan'
